# Rename the "dateCreated" column header to "date" on the manifest sheet,
# and update the active selection/scroll position to match (user clicked
# on E1 after editing D1, with the view scrolled back to the left edge).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Header row lives in row 1; column D holds the "dateCreated" label.
$ws.Range("D1").Value = "date"

# Reset the view: scroll fully left and select E1 (next to the renamed
# column), replacing the previous P1:P1048576 selection scrolled to G1.
$win = $excel.ActiveWindow
$win.ScrollColumn = 1
$win.ScrollRow = 1
[void]$ws.Range("E1").Select()
